$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Energy Efficiency in industry" / "IND_EE" row to
# "Energy Efficiency in commercial" / "COM_EE"
$ws.Range("A4").Value = "Energy Efficiency in commercial"
$ws.Range("B4").Value = "COM_EE"

# Update the active cell/selection to B5
$ws.Range("B5").Select()
